# Append rows 165-172 to the "dataset" sheet (communitySmellsDataset),
# extending the tensorflow/ranking entries per the updated web-service date handling.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the formatting already applied to column A of the existing data rows
# (bold/centered/bordered "id" style) for the new id cells.
$ws.Range("A164").Copy()
$ws.Range("A165:A172").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A165").Value = 164
$ws.Range("B165").Value = "https://github.com/tensorflow/ranking"
$ws.Range("C165").Value = "ranking"
$ws.Range("D165").Value = "tensorflow"
$ws.Range("E165").Value = "'12/03/2018"
$ws.Range("F165").Value = "'0"
$ws.Range("G165").Value = "'0"
$ws.Range("H165").Value = "'0"
$ws.Range("I165").Value = "'1"
$ws.Range("J165").Value = "'1"
$ws.Range("K165").Value = "'1"
$ws.Range("L165").Value = "'0"
$ws.Range("M165").Value = "'0"
$ws.Range("N165").Value = "'1"
$ws.Range("O165").Value = "'0"

$ws.Range("A166").Value = 165
$ws.Range("B166").Value = "https://github.com/tensorflow/ranking"
$ws.Range("C166").Value = "ranking"
$ws.Range("D166").Value = "tensorflow"
$ws.Range("E166").Value = "'01/12/2020"
$ws.Range("F166").Value = "'0"
$ws.Range("G166").Value = "'0"
$ws.Range("H166").Value = "'0"
$ws.Range("I166").Value = "'1"
$ws.Range("J166").Value = "'0"
$ws.Range("K166").Value = "'0"
$ws.Range("L166").Value = "'0"
$ws.Range("M166").Value = "'0"
$ws.Range("N166").Value = "'1"
$ws.Range("O166").Value = "'0"

$ws.Range("A167").Value = 166
$ws.Range("B167").Value = "https://github.com/tensorflow/ranking"
$ws.Range("C167").Value = "ranking"
$ws.Range("D167").Value = "tensorflow"
$ws.Range("E167").Value = "'12/03/2018"
$ws.Range("F167").Value = "'0"
$ws.Range("G167").Value = "'0"
$ws.Range("H167").Value = "'0"
$ws.Range("I167").Value = "'1"
$ws.Range("J167").Value = "'1"
$ws.Range("K167").Value = "'1"
$ws.Range("L167").Value = "'0"
$ws.Range("M167").Value = "'0"
$ws.Range("N167").Value = "'1"
$ws.Range("O167").Value = "'0"

$ws.Range("A168").Value = 167
$ws.Range("B168").Value = "https://github.com/tensorflow/ranking"
$ws.Range("C168").Value = "ranking"
$ws.Range("D168").Value = "tensorflow"
$ws.Range("E168").Value = "'12/03/2018"
$ws.Range("F168").Value = "'0"
$ws.Range("G168").Value = "'0"
$ws.Range("H168").Value = "'0"
$ws.Range("I168").Value = "'1"
$ws.Range("J168").Value = "'1"
$ws.Range("K168").Value = "'1"
$ws.Range("L168").Value = "'0"
$ws.Range("M168").Value = "'0"
$ws.Range("N168").Value = "'1"
$ws.Range("O168").Value = "'0"

$ws.Range("A169").Value = 168
$ws.Range("B169").Value = "https://github.com/tensorflow/ranking"
$ws.Range("C169").Value = "ranking"
$ws.Range("D169").Value = "tensorflow"
$ws.Range("E169").Value = "'01/12/2020"
$ws.Range("F169").Value = "'0"
$ws.Range("G169").Value = "'0"
$ws.Range("H169").Value = "'0"
$ws.Range("I169").Value = "'1"
$ws.Range("J169").Value = "'0"
$ws.Range("K169").Value = "'0"
$ws.Range("L169").Value = "'0"
$ws.Range("M169").Value = "'0"
$ws.Range("N169").Value = "'1"
$ws.Range("O169").Value = "'0"

$ws.Range("A170").Value = 169
$ws.Range("B170").Value = "https://github.com/tensorflow/ranking"
$ws.Range("C170").Value = "ranking"
$ws.Range("D170").Value = "tensorflow"
$ws.Range("E170").Value = "'01/12/2020"
$ws.Range("F170").Value = "'0"
$ws.Range("G170").Value = "'0"
$ws.Range("H170").Value = "'0"
$ws.Range("I170").Value = "'1"
$ws.Range("J170").Value = "'0"
$ws.Range("K170").Value = "'0"
$ws.Range("L170").Value = "'0"
$ws.Range("M170").Value = "'0"
$ws.Range("N170").Value = "'1"
$ws.Range("O170").Value = "'0"

$ws.Range("A171").Value = 170
$ws.Range("B171").Value = "https://github.com/tensorflow/ranking"
$ws.Range("C171").Value = "ranking"
$ws.Range("D171").Value = "tensorflow"
$ws.Range("E171").Value = "'01/12/2020"
$ws.Range("F171").Value = "'0"
$ws.Range("G171").Value = "'0"
$ws.Range("H171").Value = "'0"
$ws.Range("I171").Value = "'1"
$ws.Range("J171").Value = "'0"
$ws.Range("K171").Value = "'0"
$ws.Range("L171").Value = "'0"
$ws.Range("M171").Value = "'0"
$ws.Range("N171").Value = "'1"
$ws.Range("O171").Value = "'0"

$ws.Range("A172").Value = 171
$ws.Range("B172").Value = "https://github.com/tensorflow/ranking"
$ws.Range("C172").Value = "ranking"
$ws.Range("D172").Value = "tensorflow"
$ws.Range("E172").Value = "'12/03/2018"
$ws.Range("F172").Value = "'0"
$ws.Range("G172").Value = "'0"
$ws.Range("H172").Value = "'0"
$ws.Range("I172").Value = "'1"
$ws.Range("J172").Value = "'1"
$ws.Range("K172").Value = "'1"
$ws.Range("L172").Value = "'0"
$ws.Range("M172").Value = "'0"
$ws.Range("N172").Value = "'1"
$ws.Range("O172").Value = "'0"
